$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13, shifting everything below down by one
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the bibo:status metadata
$ws.Cells.Item(13, 1).Value = "http://purl.org/ontology/bibo/status"
$ws.Cells.Item(13, 2).Value = "Draft Controlled Vocabulary"

# Copy the style of column A label cells (e.g. A12) onto the new A13 cell
$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The row-insert operation copies formatting down from row 12 into the
# new row 13 (B13 picks up the hyperlink style, C13 becomes populated
# with an empty styled cell). Clear that unwanted formatting/content so
# B13 is a plain cell and C13 stays empty.
$ws.Cells.Item(13, 2).ClearFormats()
$ws.Cells.Item(13, 3).Clear()

# Update the selection to match the target state
$ws.Range("A13:XFD13").Select()
